$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2020 column (N) of data, mirroring the formatting of column M
# (copy format only, then set the value) for every row that already has
# data in column M.

function Set-NCell($row, $value) {
    $ws.Range("M$row").Copy() | Out-Null
    $ws.Range("N$row").PasteSpecial(-4122) | Out-Null
    if ($null -ne $value) {
        $ws.Range("N$row").Value = $value
    }
}

Set-NCell 3  2020
Set-NCell 5  2198.6999999999998
Set-NCell 6  132.69999999999999
Set-NCell 7  242.9
Set-NCell 8  203.3
Set-NCell 9  202.8
Set-NCell 10 284.7
Set-NCell 11 294.89999999999998
Set-NCell 12 802.5
Set-NCell 13 28.1
Set-NCell 14 6.8
Set-NCell 15 $null
Set-NCell 16 27.4
Set-NCell 17 17.5
Set-NCell 18 24.7
Set-NCell 19 31.5
Set-NCell 20 30.4
Set-NCell 21 24.8
Set-NCell 22 30.7
Set-NCell 23 30.1
Set-NCell 24 21.2
Set-NCell 25 11.6

# Match the author's final selection in the saved workbook.
$ws.Range("M25").Select() | Out-Null
